$wb = $excel.ActiveWorkbook

# --- Sheet 1 (汽車 / Car) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B1").Value = "中華FB308W"
$ws1.Range("C1").Value = 2835
$ws1.Range("D1").Value = "陳亭妃"
$ws1.Range("E1").Value = "98年03月24日"
$ws1.Range("F1").Value = "繼承"
$ws1.Range("G1").Value = 10000
$ws1.Rows.Item(2).Delete()

# --- Sheet 2 (債務 / Debt) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B1").Value = "中期放款"
$ws2.Range("C1").Value = "陳亭妃"
$ws2.Range("D1").Value = "合作金庫商業銀行臺南市北區西門路"
$ws2.Range("E1").Value = 970000
$ws2.Range("F1").Value = "89年03月29日"
$ws2.Range("G1").Value = "信用貸款"
$ws2.Rows.Item(2).Delete()
